# Scheduled GitHub Actions symbol-list refresh (Sat Dec 17 10:53:26 UTC 2022):
# updates the cached coin price snapshots (column D) and a couple of
# "<rank><name><symbol>[Worstin24h]" labels (column E) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $value) {
    # These columns store plain text look-alikes of numbers (e.g. "236.62"),
    # so force text formatting while writing, then drop back to the default
    # "Normal" style so no stray per-cell formatting is left behind.
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "236.62"
Set-TextValue "D3" "21.88"
Set-TextValue "D4" "5.454"
Set-TextValue "D5" "0.05634"
Set-TextValue "D7" "3.345"
Set-TextValue "D8" "1.071"
Set-TextValue "D10" "0.1395"
Set-TextValue "D11" "0.07341"
Set-TextValue "D13" "0.02972"
Set-TextValue "D14" "0.09239"
Set-TextValue "D15" "0.001669"
Set-TextValue "D16" "3.252"
Set-TextValue "D17" "0.04760"
Set-TextValue "D18" "0.0005745"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006227"
Set-TextValue "D20" "0.005091"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "3.873"
Set-TextValue "D26" "0.1054"
Set-TextValue "D27" "0.0004013"
Set-TextValue "D40" "0.04119"
Set-TextValue "D41" "0.006979"
Set-TextValue "D42" "0.003504"
Set-TextValue "D43" "0.1036"
Set-TextValue "D44" "0.008911"
Set-TextValue "D45" "0.00005441"
Set-TextValue "D47" "0.6757"
Set-TextValue "D48" "0.03779"
Set-TextValue "E48" "47BOLOBOLO"
Set-TextValue "D49" "0.00002102"
